$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-number-looking price strings to stay text (matches source inlineStr cells)
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '64.244.98'
$ws.Range('E2').Value = '  +0.95%  '
$ws.Range('D3').Value = '3.485.74'
$ws.Range('E3').Value = '  +0.37%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '585.92'
$ws.Range('E5').Value = '  +0.56%  '
$ws.Range('D6').Value = '134.01'
$ws.Range('E6').Value = '  +2.15%  '
$ws.Range('D7').Value = '3.484.89'
$ws.Range('E7').Value = '  +0.34%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '0.485'
$ws.Range('E9').Value = '  -0.60%  '
$ws.Range('E10').Value = '  +0.24%  '
$ws.Range('D11').Value = '7.19'
$ws.Range('E11').Value = '  +2.28%  '
$ws.Range('E12').Value = '  -2.21%  '
$ws.Range('D13').Value = '4.079.79'
$ws.Range('E13').Value = '  +0.57%  '
$ws.Range('E14').Value = '  +2.20%  '
$ws.Range('E15').Value = '  +1.25%  '
$ws.Range('D16').Value = '3.484.40'
$ws.Range('E16').Value = '  +0.54%  '
$ws.Range('D17').Value = '64.305.35'
$ws.Range('E17').Value = '  +0.83%  '
$ws.Range('D18').Value = '25.13'
$ws.Range('E18').Value = '  -9.19%  '
$ws.Range('D19').Value = '10.00'
$ws.Range('E19').Value = '  +1.07%  '
$ws.Range('D20').Value = '5.69'
$ws.Range('E20').Value = '  +0.85%  '
$ws.Range('E21').Value = '  -4.17%  '
$ws.Range('D22').Value = '384.02'
$ws.Range('E22').Value = '  -1.14%  '
$ws.Range('E23').Value = '  -1.58%  '
$ws.Range('D24').Value = '3.625.07'
$ws.Range('E24').Value = '  +0.46%  '
$ws.Range('D25').Value = '74.11'
$ws.Range('E25').Value = '  +1.97%  '
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('E27').Value = '  -0.72%  '
$ws.Range('D28').Value = '0.0000112'
$ws.Range('E28').Value = '  +3.56%  '
$ws.Range('D29').Value = '1.55'
$ws.Range('E29').Value = '  -0.25%  '
$ws.Range('E30').Value = '  +0.17%  '
$ws.Range('D31').Value = '7.42'
$ws.Range('E31').Value = '  +0.53%  '
$ws.Range('E32').Value = '  -0.29%  '
$ws.Range('D33').Value = '8.20'
$ws.Range('E33').Value = '  +0.75%  '
$ws.Range('D34').Value = '3.507.76'
$ws.Range('E34').Value = '  +1.04%  '
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('E36').Value = '  +2.44%  '
$ws.Range('D37').Value = '23.37'
$ws.Range('E37').Value = '  -1.11%  '
$ws.Range('D38').Value = '5.26'
$ws.Range('E38').Value = '  -0.30%  '
$ws.Range('D39').Value = '6.85'
$ws.Range('E39').Value = '  -1.46%  '
$ws.Range('E40').Value = '  -1.73%  '
$ws.Range('D41').Value = '162.13'
$ws.Range('E41').Value = '  -4.11%  '
$ws.Range('D42').Value = '0.0778'
$ws.Range('E42').Value = '  -2.92%  '
$ws.Range('D43').Value = '0.804'
$ws.Range('E43').Value = '  -0.35%  '
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('D45').Value = '25.42'
$ws.Range('E45').Value = '  -1.11%  '
$ws.Range('E46').Value = '  +0.27%  '
$ws.Range('D47').Value = '4.39'
$ws.Range('E47').Value = '  +1.29%  '
$ws.Range('E48').Value = '  +0.45%  '
$ws.Range('E49').Value = '  +0.90%  '
$ws.Range('D50').Value = '2.466.35'
$ws.Range('E50').Value = '  +2.33%  '
$ws.Range('D51').Value = '6.72'
$ws.Range('E51').Value = '  -1.88%  '
